# Edit script: populate MayRaw sheet with raw May 2023 statistics data.
# This mirrors the data already present in other *Raw sheets (JanuaryRaw, AprilRaw, etc).
# Downstream sheets (May, Yearly total) use formulas referencing MayRaw and will
# recalculate automatically once these values are populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MayRaw")

# Header row (row 1) - labels already present as shared strings in the workbook
$ws.Range("A1").Value = "Library"
$ws.Range("B1").Value = "Items owned by this library checked out at this library this month"
$ws.Range("C1").Value = "Items owned by other libraries checked out at this library this month"
$ws.Range("D1").Value = "Total circulation this month"

# Data rows
$ws.Range("A2").Value = 'Atchison Public Library'
$ws.Range("B2").Value = 4944
$ws.Range("C2").Value = 1613
$ws.Range("D2").Value = 6557

$ws.Range("A3").Value = 'Baldwin City Public Library'
$ws.Range("B3").Value = 2901
$ws.Range("C3").Value = 476
$ws.Range("D3").Value = 3377

$ws.Range("A4").Value = 'Basehor Community Library'
$ws.Range("B4").Value = 9600
$ws.Range("C4").Value = 1279
$ws.Range("D4").Value = 10879

$ws.Range("A5").Value = 'Bern Community Library'
$ws.Range("B5").Value = 121
$ws.Range("C5").Value = 132
$ws.Range("D5").Value = 253

$ws.Range("A6").Value = 'Bonner Springs City Library'
$ws.Range("B6").Value = 5424
$ws.Range("C6").Value = 1218
$ws.Range("D6").Value = 6642

$ws.Range("A7").Value = 'Burlingame Community Library'
$ws.Range("B7").Value = 511
$ws.Range("C7").Value = 195
$ws.Range("D7").Value = 706

$ws.Range("A8").Value = 'Carbondale City Library'
$ws.Range("B8").Value = 672
$ws.Range("C8").Value = 122
$ws.Range("D8").Value = 794

$ws.Range("A9").Value = 'Centralia Community Library'
$ws.Range("B9").Value = 344
$ws.Range("C9").Value = 41
$ws.Range("D9").Value = 385

$ws.Range("A10").Value = 'Corning City Library'
$ws.Range("B10").Value = 20
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 25

$ws.Range("A11").Value = 'Digital Content'

$ws.Range("A12").Value = 'Doniphan County Library - Elwood'
$ws.Range("B12").Value = 45
$ws.Range("C12").Value = 8
$ws.Range("D12").Value = 53

$ws.Range("A13").Value = 'Doniphan County Library - Highland'
$ws.Range("B13").Value = 261
$ws.Range("C13").Value = 171
$ws.Range("D13").Value = 432

$ws.Range("A14").Value = 'Doniphan County Library - Troy'
$ws.Range("B14").Value = 535
$ws.Range("C14").Value = 111
$ws.Range("D14").Value = 646

$ws.Range("A15").Value = 'Doniphan County Library - Wathena'
$ws.Range("B15").Value = 343
$ws.Range("C15").Value = 77
$ws.Range("D15").Value = 420

$ws.Range("A16").Value = 'Effingham Community Library'
$ws.Range("B16").Value = 425
$ws.Range("C16").Value = 61
$ws.Range("D16").Value = 486

$ws.Range("A17").Value = 'Eudora Community Library'
$ws.Range("B17").Value = 2013
$ws.Range("C17").Value = 555
$ws.Range("D17").Value = 2568

$ws.Range("A18").Value = 'Everest, Barnes Reading Room'
$ws.Range("B18").Value = 109
$ws.Range("C18").Value = 89
$ws.Range("D18").Value = 198

$ws.Range("A19").Value = 'Hiawatha, Morrill Public Library'
$ws.Range("B19").Value = 1808
$ws.Range("C19").Value = 707
$ws.Range("D19").Value = 2515

$ws.Range("A20").Value = 'Highland Community College'
$ws.Range("B20").Value = 9
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 12

$ws.Range("A21").Value = 'Holton, Beck-Bookman Library'
$ws.Range("B21").Value = 1880
$ws.Range("C21").Value = 496
$ws.Range("D21").Value = 2376

$ws.Range("A22").Value = 'Horton Public Library'
$ws.Range("B22").Value = 99
$ws.Range("C22").Value = 32
$ws.Range("D22").Value = 131

$ws.Range("A23").Value = 'Lansing Community Library'
$ws.Range("B23").Value = 1866
$ws.Range("C23").Value = 763
$ws.Range("D23").Value = 2629

$ws.Range("A24").Value = 'Leavenworth Public Library'
$ws.Range("B24").Value = 8722
$ws.Range("C24").Value = 1907
$ws.Range("D24").Value = 10629

$ws.Range("A25").Value = 'Linwood Community Library'
$ws.Range("B25").Value = 665
$ws.Range("C25").Value = 176
$ws.Range("D25").Value = 841

$ws.Range("A26").Value = 'Louisburg Library'

$ws.Range("A27").Value = 'Lyndon Carnegie Library'
$ws.Range("B27").Value = 424
$ws.Range("C27").Value = 209
$ws.Range("D27").Value = 633

$ws.Range("A28").Value = 'McLouth Public Library'
$ws.Range("B28").Value = 168
$ws.Range("C28").Value = 53
$ws.Range("D28").Value = 221

$ws.Range("A29").Value = 'Meriden-Ozawkie Public Library'
$ws.Range("B29").Value = 1785
$ws.Range("C29").Value = 449
$ws.Range("D29").Value = 2234

$ws.Range("A30").Value = 'Northeast Kansas Library System'
$ws.Range("B30").Value = 17
$ws.Range("C30").Value = 18
$ws.Range("D30").Value = 35

$ws.Range("A31").Value = 'Nortonville Public Library'
$ws.Range("B31").Value = 242
$ws.Range("C31").Value = 61
$ws.Range("D31").Value = 303

$ws.Range("A32").Value = 'Osage City Library'
$ws.Range("B32").Value = 2334
$ws.Range("C32").Value = 440
$ws.Range("D32").Value = 2774

$ws.Range("A33").Value = 'Osawatomie Public Library'
$ws.Range("B33").Value = 908
$ws.Range("C33").Value = 363
$ws.Range("D33").Value = 1271

$ws.Range("A34").Value = 'Oskaloosa Public Library'
$ws.Range("B34").Value = 624
$ws.Range("C34").Value = 171
$ws.Range("D34").Value = 795

$ws.Range("A35").Value = 'Ottawa Library'
$ws.Range("B35").Value = 6667
$ws.Range("C35").Value = 880
$ws.Range("D35").Value = 7547

$ws.Range("A36").Value = 'Overbrook Public Library'
$ws.Range("B36").Value = 980
$ws.Range("C36").Value = 241
$ws.Range("D36").Value = 1221

$ws.Range("A37").Value = 'Paola Free Library'
$ws.Range("B37").Value = 3138
$ws.Range("C37").Value = 467
$ws.Range("D37").Value = 3605

$ws.Range("A38").Value = 'Perry-Lecompton Community Library'
$ws.Range("B38").Value = 50
$ws.Range("C38").Value = 11
$ws.Range("D38").Value = 61

$ws.Range("A39").Value = 'Pomona Community Library'
$ws.Range("B39").Value = 206
$ws.Range("C39").Value = 51
$ws.Range("D39").Value = 257

$ws.Range("A40").Value = 'Prairie Hills Schools - Axtell Public School'
$ws.Range("B40").Value = 115
$ws.Range("C40").Value = 1
$ws.Range("D40").Value = 116

$ws.Range("A41").Value = 'Prairie Hills Schools - Sabetha Elementary School'
$ws.Range("B41").Value = 405
$ws.Range("C41").Value = 15
$ws.Range("D41").Value = 420

$ws.Range("A42").Value = 'Prairie Hills Schools - Sabetha High School'
$ws.Range("B42").Value = 14
$ws.Range("C42").Value = 3
$ws.Range("D42").Value = 17

$ws.Range("A43").Value = 'Prairie Hills Schools - Sabetha Middle School'
$ws.Range("B43").Value = 2
$ws.Range("D43").Value = 2

$ws.Range("A44").Value = 'Prairie Hills Schools - Wetmore Academic Center (Permanently closed)'

$ws.Range("A45").Value = 'Richmond Public Library'
$ws.Range("B45").Value = 465
$ws.Range("C45").Value = 77
$ws.Range("D45").Value = 542

$ws.Range("A46").Value = 'Rossville Community Library'
$ws.Range("B46").Value = 1177
$ws.Range("C46").Value = 401
$ws.Range("D46").Value = 1578

$ws.Range("A47").Value = 'Sabetha, Mary Cotton Library'
$ws.Range("B47").Value = 4403
$ws.Range("C47").Value = 959
$ws.Range("D47").Value = 5362

$ws.Range("A48").Value = 'Seneca Free Library'
$ws.Range("B48").Value = 2300
$ws.Range("C48").Value = 231
$ws.Range("D48").Value = 2531

$ws.Range("A49").Value = 'Silver Lake Library'
$ws.Range("B49").Value = 1168
$ws.Range("C49").Value = 311
$ws.Range("D49").Value = 1479

$ws.Range("A50").Value = 'Tonganoxie Public Library'
$ws.Range("B50").Value = 4273
$ws.Range("C50").Value = 815
$ws.Range("D50").Value = 5088

$ws.Range("A51").Value = 'Valley Falls, Delaware Township Library'
$ws.Range("B51").Value = 465
$ws.Range("C51").Value = 293
$ws.Range("D51").Value = 758

$ws.Range("A52").Value = 'Wellsville City Library'
$ws.Range("B52").Value = 1790
$ws.Range("C52").Value = 475
$ws.Range("D52").Value = 2265

$ws.Range("A53").Value = 'Wetmore Public Library'
$ws.Range("B53").Value = 103
$ws.Range("C53").Value = 136
$ws.Range("D53").Value = 239

$ws.Range("A54").Value = 'Williamsburg Community Library'
$ws.Range("B54").Value = 239
$ws.Range("C54").Value = 22
$ws.Range("D54").Value = 261

$ws.Range("A55").Value = 'Winchester Public Library'
$ws.Range("B55").Value = 447
$ws.Range("C55").Value = 278
$ws.Range("D55").Value = 725

# Force recalculation so dependent sheets (May, Yearly total) pick up the new values
$excel.Calculate()

